$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E (shifts old E:I -> F:J), this brings in the new
# "slrtype" column between "Study_Types" (D) and "slrtype_Radio_button" (now F).
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = 17

# --- Row 1 (headers) ---
$ws.Range("E1").Value = "slrtype"
# The freshly inserted column inherits D's bold/centered header style; the
# target file keeps this header unstyled (matching its neighbour F1).
$ws.Range("E1").Style = "Normal"

# --- Column A: population ids ---
$ws.Range("A3").Value = "pop2"
$ws.Range("A4").Value = "pop3"
$ws.Range("A5").Value = "pop4"

# --- Columns B & C: fill down the same values used in row 2 ---
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = $ws.Range("C2").Text
$ws.Range("B4").Value = $ws.Range("B2").Text
$ws.Range("C4").Value = $ws.Range("C2").Text
$ws.Range("B5").Value = $ws.Range("B2").Text
$ws.Range("C5").Value = $ws.Range("C2").Text

# --- New column E: mirror the Study_Types values from column D ---
$ws.Range("E2").Value = $ws.Range("D2").Text
$ws.Range("E3").Value = $ws.Range("D3").Text
$ws.Range("E4").Value = "Quality of Life"
$ws.Range("E5").Value = $ws.Range("D5").Text

# --- Column H (previously G) "Files_to_upload": swap Cochrane paths for the
#     new SearchStrategy template paths for QOL and RWE ---
$ws.Range("H4").Value = "\Testdata\Templates\SearchStrategy\QOL.xlsx"
$ws.Range("H5").Value = "\Testdata\Templates\SearchStrategy\RWE.xlsx"

# --- Update the view: drop the old scroll/selection and select C7 ---
$ws.Range("C7").Select()
